# Update NATMI ligand-receptor edge statistics for C3-Cd46 sheet
# following recalculation with updated cell-expression counts (Dr Hou advice).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.942832
$ws.Range("H2").Value = 122.828496
$ws.Range("I2").Value = 0.2583000005785167
$ws.Range("J2").Value = 0.2583000005785167
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.452892666666665
$ws.Range("N2").Value = 25.358678
$ws.Range("O2").Value = 0.5664982795292011
$ws.Range("P2").Value = 0.566498279529201
$ws.Range("Q2").Value = 346.0853643653653
$ws.Range("R2").Value = 3114.768279288288
$ws.Range("S2").Value = 0.1463265059301213
$ws.Range("T2").Value = 0.1463265059301213

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.942832
$ws.Range("H3").Value = 122.828496
$ws.Range("I3").Value = 0.2583000005785167
$ws.Range("J3").Value = 0.2583000005785167
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.429517666666667
$ws.Range("N3").Value = 10.288553
$ws.Range("O3").Value = 0.2298403557687432
$ws.Range("P3").Value = 0.2298403557687431
$ws.Range("Q3").Value = 140.4141656673654
$ws.Range("R3").Value = 1263.727491006288
$ws.Range("S3").Value = 0.05936776402803284
$ws.Range("T3").Value = 0.05936776402803283

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.942832
$ws.Range("H4").Value = 122.828496
$ws.Range("I4").Value = 0.2583000005785167
$ws.Range("J4").Value = 0.2583000005785167
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.6234873333333334
$ws.Range("N4").Value = 1.870462
$ws.Range("O4").Value = 0.04178504514015868
$ws.Range("P4").Value = 0.04178504514015867
$ws.Range("Q4").Value = 25.52733714279467
$ws.Range("R4").Value = 229.746034285152
$ws.Range("S4").Value = 0.01079307718387633
$ws.Range("T4").Value = 0.01079307718387633

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.942832
$ws.Range("H5").Value = 122.828496
$ws.Range("I5").Value = 0.2583000005785167
$ws.Range("J5").Value = 0.2583000005785167
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.415405666666667
$ws.Range("N5").Value = 7.246217000000001
$ws.Range("O5").Value = 0.1618763195618971
$ws.Range("P5").Value = 0.1618763195618971
$ws.Range("Q5").Value = 98.89354842218135
$ws.Range("R5").Value = 890.041935799632
$ws.Range("S5").Value = 0.04181265343648616
$ws.Range("T5").Value = 0.04181265343648616

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 79.68771233333334
$ws.Range("H6").Value = 239.063137
$ws.Range("I6").Value = 0.5027335710876245
$ws.Range("J6").Value = 0.5027335710876245
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.452892666666665
$ws.Range("N6").Value = 25.358678
$ws.Range("O6").Value = 0.5664982795292011
$ws.Range("P6").Value = 0.566498279529201
$ws.Range("Q6").Value = 673.5916792058762
$ws.Range("R6").Value = 6062.325112852885
$ws.Range("S6").Value = 0.2847977030827106
$ws.Range("T6").Value = 0.2847977030827105

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 79.68771233333334
$ws.Range("H7").Value = 239.063137
$ws.Range("I7").Value = 0.5027335710876245
$ws.Range("J7").Value = 0.5027335710876245
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.429517666666667
$ws.Range("N7").Value = 10.288553
$ws.Range("O7").Value = 0.2298403557687432
$ws.Range("P7").Value = 0.2298403557687431
$ws.Range("Q7").Value = 273.2904172634179
$ws.Range("R7").Value = 2459.613755370761
$ws.Range("S7").Value = 0.1155484628356704
$ws.Range("T7").Value = 0.1155484628356703

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 79.68771233333334
$ws.Range("H8").Value = 239.063137
$ws.Range("I8").Value = 0.5027335710876245
$ws.Range("J8").Value = 0.5027335710876245
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.6234873333333334
$ws.Range("N8").Value = 1.870462
$ws.Range("O8").Value = 0.04178504514015868
$ws.Range("P8").Value = 0.04178504514015867
$ws.Range("Q8").Value = 49.68427926214378
$ws.Range("R8").Value = 447.158513359294
$ws.Range("S8").Value = 0.02100674496136956
$ws.Range("T8").Value = 0.02100674496136956

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 79.68771233333334
$ws.Range("H9").Value = 239.063137
$ws.Range("I9").Value = 0.5027335710876245
$ws.Range("J9").Value = 0.5027335710876245
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.415405666666667
$ws.Range("N9").Value = 7.246217000000001
$ws.Range("O9").Value = 0.1618763195618971
$ws.Range("P9").Value = 0.1618763195618971
$ws.Range("Q9").Value = 192.4781519336366
$ws.Range("R9").Value = 1732.303367402729
$ws.Range("S9").Value = 0.08138066020787402
$ws.Range("T9").Value = 0.081380660207874

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.53186833333334
$ws.Range("H10").Value = 112.595605
$ws.Range("I10").Value = 0.2367809244903433
$ws.Range("J10").Value = 0.2367809244903433
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.452892666666665
$ws.Range("N10").Value = 25.358678
$ws.Range("O10").Value = 0.5664982795292011
$ws.Range("P10").Value = 0.566498279529201
$ws.Range("Q10").Value = 317.2528546011322
$ws.Range("R10").Value = 2855.27569141019
$ws.Range("S10").Value = 0.1341359863491131
$ws.Range("T10").Value = 0.1341359863491131

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 37.53186833333334
$ws.Range("H11").Value = 112.595605
$ws.Range("I11").Value = 0.2367809244903433
$ws.Range("J11").Value = 0.2367809244903433
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.429517666666667
$ws.Range("N11").Value = 10.288553
$ws.Range("O11").Value = 0.2298403557687432
$ws.Range("P11").Value = 0.2298403557687431
$ws.Range("Q11").Value = 128.7162055121739
$ws.Range("R11").Value = 1158.445849609565
$ws.Range("S11").Value = 0.05442181192411243
$ws.Range("T11").Value = 0.05442181192411241

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 37.53186833333334
$ws.Range("H12").Value = 112.595605
$ws.Range("I12").Value = 0.2367809244903433
$ws.Range("J12").Value = 0.2367809244903433
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.6234873333333334
$ws.Range("N12").Value = 1.870462
$ws.Range("O12").Value = 0.04178504514015868
$ws.Range("P12").Value = 0.04178504514015867
$ws.Range("Q12").Value = 23.40064450216778
$ws.Range("R12").Value = 210.60580051951
$ws.Range("S12").Value = 0.009893901618157497
$ws.Range("T12").Value = 0.009893901618157494

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 37.53186833333334
$ws.Range("H13").Value = 112.595605
$ws.Range("I13").Value = 0.2367809244903433
$ws.Range("J13").Value = 0.2367809244903433
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.415405666666667
$ws.Range("N13").Value = 7.246217000000001
$ws.Range("O13").Value = 0.1618763195618971
$ws.Range("P13").Value = 0.1618763195618971
$ws.Range("Q13").Value = 90.65468745292057
$ws.Range("R13").Value = 815.8921870762852
$ws.Range("S13").Value = 0.03832922459896024
$ws.Range("T13").Value = 0.03832922459896023

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.3464216666666666
$ws.Range("H14").Value = 1.039265
$ws.Range("I14").Value = 0.002185503843515531
$ws.Range("J14").Value = 0.002185503843515531
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.452892666666665
$ws.Range("N14").Value = 25.358678
$ws.Range("O14").Value = 0.5664982795292011
$ws.Range("P14").Value = 0.566498279529201
$ws.Range("Q14").Value = 2.92826516574111
$ws.Range("R14").Value = 26.35438649167
$ws.Range("S14").Value = 0.001238084167256005
$ws.Range("T14").Value = 0.001238084167256005

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.3464216666666666
$ws.Range("H15").Value = 1.039265
$ws.Range("I15").Value = 0.002185503843515531
$ws.Range("J15").Value = 0.002185503843515531
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.429517666666667
$ws.Range("N15").Value = 10.288553
$ws.Range("O15").Value = 0.2298403557687432
$ws.Range("P15").Value = 0.2298403557687431
$ws.Range("Q15").Value = 1.188059225949444
$ws.Range("R15").Value = 10.692533033545
$ws.Range("S15").Value = 0.0005023169809275655
$ws.Range("T15").Value = 0.0005023169809275652

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.3464216666666666
$ws.Range("H16").Value = 1.039265
$ws.Range("I16").Value = 0.002185503843515531
$ws.Range("J16").Value = 0.002185503843515531
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.6234873333333334
$ws.Range("N16").Value = 1.870462
$ws.Range("O16").Value = 0.04178504514015868
$ws.Range("P16").Value = 0.04178504514015867
$ws.Range("Q16").Value = 0.2159895211588889
$ws.Range("R16").Value = 1.94390569043
$ws.Range("S16").Value = 0.00009132137675528676
$ws.Range("T16").Value = 0.00009132137675528675

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.3464216666666666
$ws.Range("H17").Value = 1.039265
$ws.Range("I17").Value = 0.002185503843515531
$ws.Range("J17").Value = 0.002185503843515531
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.415405666666667
$ws.Range("N17").Value = 7.246217000000001
$ws.Range("O17").Value = 0.1618763195618971
$ws.Range("P17").Value = 0.1618763195618971
$ws.Range("Q17").Value = 0.8367488567227778
$ws.Range("R17").Value = 7.530739710504999
$ws.Range("S17").Value = 0.0003537813185766745
$ws.Range("T17").Value = 0.0003537813185766744
